$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H58").Value = 1734
$ws.Range("I58").Value = 297.33334
$ws.Range("J58").Value = 4320
$ws.Range("K58").Value = 892.0000200000001
$ws.Range("L58").Value = 12960
$ws.Range("M58").Value = -742.0000200000001
$ws.Range("N58").Value = -13260

$ws.Range("H74").Value = 4198.2
$ws.Range("I74").Value = 4198
$ws.Range("J74").Value = 4198.5
$ws.Range("K74").Value = 4198
$ws.Range("L74").Value = 4198.5
$ws.Range("M74").Value = -3262
$ws.Range("N74").Value = -6070.5

$ws.Range("H77").Value = 4198.2
$ws.Range("I77").Value = 4198
$ws.Range("J77").Value = 4198.5
$ws.Range("K77").Value = 20990
$ws.Range("L77").Value = 20992.5
$ws.Range("M77").Value = -16310
$ws.Range("N77").Value = -30352.5

$ws.Range("H93").Value = 84600
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 84600
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 84600
$ws.Range("N93").Value = -89592

$ws.Range("H137").Value = 1808.25
$ws.Range("I137").Value = 1305.5
$ws.Range("J137").Value = 1975.8334
$ws.Range("K137").Value = 3916.5
$ws.Range("L137").Value = 5927.5002
$ws.Range("M137").Value = -1366.5

$ws.Range("H138").Value = 3265.04
$ws.Range("I138").Value = 4007.1333
$ws.Range("J138").Value = 2151.9
$ws.Range("K138").Value = 12021.3999
$ws.Range("L138").Value = 6455.700000000001
$ws.Range("M138").Value = -6881.3999
$ws.Range("N138").Value = -16735.7

$ws.Range("H139").Value = 45000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 45000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 45000
$ws.Range("N139").Value = -55280

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H61").Value = 4758
$ws.Range("I61").Value = 3641.4167
$ws.Range("J61").Value = 7437.8
$ws.Range("K61").Value = 3641.4167
$ws.Range("L61").Value = 7437.8
$ws.Range("M61").Value = -3429.4167
$ws.Range("N61").Value = -7861.8

$ws.Range("H74").Value = 3330
$ws.Range("I74").Value = 3224.158
$ws.Range("J74").Value = 4000.3333
$ws.Range("K74").Value = 3224.158
$ws.Range("L74").Value = 4000.3333
$ws.Range("M74").Value = -2350.158

$ws.Range("H77").Value = 3330
$ws.Range("I77").Value = 3224.158
$ws.Range("J77").Value = 4000.3333
$ws.Range("K77").Value = 16120.79
$ws.Range("L77").Value = 20001.6665
$ws.Range("M77").Value = -11752.79

$ws.Range("H88").Value = 2774.1875
$ws.Range("I88").Value = 1905.7142
$ws.Range("J88").Value = 3449.6667
$ws.Range("K88").Value = 1905.7142
$ws.Range("L88").Value = 3449.6667
$ws.Range("M88").Value = -1499.7142

$ws.Range("H91").Value = 2774.1875
$ws.Range("I91").Value = 1905.7142
$ws.Range("J91").Value = 3449.6667
$ws.Range("K91").Value = 1905.7142
$ws.Range("L91").Value = 3449.6667
$ws.Range("M91").Value = -501.7141999999999

$ws.Range("H132").Value = 1613.25
$ws.Range("I132").Value = 956.6667
$ws.Range("J132").Value = 3583
$ws.Range("K132").Value = 2870.0001
$ws.Range("L132").Value = 10749
$ws.Range("M132").Value = -340.0001000000002
$ws.Range("N132").Value = -15809

$ws.Range("H136").Value = 4758
$ws.Range("I136").Value = 3641.4167
$ws.Range("J136").Value = 7437.8
$ws.Range("K136").Value = 10924.2501
$ws.Range("L136").Value = 22313.4
$ws.Range("M136").Value = -8374.250100000001
$ws.Range("N136").Value = -27413.4

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("M26").ClearContents()
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0

$ws.Range("H99").Value = 3687
$ws.Range("I99").Value = 4332.6665
$ws.Range("J99").Value = 1750
$ws.Range("K99").Value = 4332.6665
$ws.Range("L99").Value = 1750
$ws.Range("M99").Value = -2834.6665

$ws.Range("H134").Value = 8637.6
$ws.Range("I134").Value = 9314.048000000001
$ws.Range("J134").Value = 5086.25
$ws.Range("K134").Value = 27942.144
$ws.Range("L134").Value = 15258.75
$ws.Range("M134").Value = -25407.144

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H22").Value = 575
$ws.Range("I22").Value = 149
$ws.Range("J22").Value = 1001
$ws.Range("K22").Value = 149
$ws.Range("L22").Value = 1001
$ws.Range("M22").Value = 201
$ws.Range("N22").Value = -1701

$ws.Range("H62").Value = 3800
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 3450
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 3450
$ws.Range("M62").Value = -3876

$ws.Range("H65").Value = 3800
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 3450
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 17250
$ws.Range("M65").Value = -19380

$ws.Range("H111").Value = 32000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 32000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 32000
$ws.Range("N111").Value = -40180

$ws.Range("H132").Value = 2415.24
$ws.Range("I132").Value = 1131.8182
$ws.Range("J132").Value = 3423.6428
$ws.Range("K132").Value = 3395.4546
$ws.Range("L132").Value = 10270.9284
$ws.Range("M132").Value = -865.4546
$ws.Range("N132").Value = -15330.9284

$ws.Range("H134").Value = 1641.5
$ws.Range("I134").Value = 1574.2307
$ws.Range("J134").Value = 1933
$ws.Range("K134").Value = 4722.6921
$ws.Range("L134").Value = 5799
$ws.Range("M134").Value = -2187.6921
$ws.Range("N134").Value = -10869

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H129").Value = 52754.5
$ws.Range("I129").Value = 952
$ws.Range("J129").Value = 73475.5
$ws.Range("K129").Value = 2856
$ws.Range("L129").Value = 220426.5
$ws.Range("M129").Value = 2144
$ws.Range("N129").Value = -230426.5

$ws.Range("H130").Value = 1976.6666
$ws.Range("I130").Value = 1430
$ws.Range("J130").Value = 2250
$ws.Range("K130").Value = 4290
$ws.Range("L130").Value = 6750
$ws.Range("M130").Value = 730
$ws.Range("N130").Value = -16790

$ws.Range("H131").Value = 11380641
$ws.Range("I131").Value = 71429144
$ws.Range("J131").Value = 20113.568
$ws.Range("K131").Value = 214287432
$ws.Range("L131").Value = 60340.704
$ws.Range("M131").Value = -214282392
$ws.Range("N131").Value = -70420.704

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("M80").ClearContents()
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3000
$ws.Range("N80").Value = -4996

$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 15000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -24984

$ws.Range("H126").Value = 58554.945
$ws.Range("I126").Value = 3185.9375
$ws.Range("J126").Value = 501507
$ws.Range("K126").Value = 9557.8125
$ws.Range("L126").Value = 1504521
$ws.Range("M126").Value = -7087.8125

$ws.Range("H132").Value = 4256.5835
$ws.Range("I132").Value = 3146.1333
$ws.Range("J132").Value = 6107.3335
$ws.Range("K132").Value = 9438.3999
$ws.Range("L132").Value = 18322.0005
$ws.Range("M132").Value = -6908.3999
$ws.Range("N132").Value = -23382.0005

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H40").Value = 6195.273
$ws.Range("I40").Value = 3390.6365
$ws.Range("J40").Value = 8999.909
$ws.Range("K40").Value = 3390.6365
$ws.Range("L40").Value = 8999.909
$ws.Range("M40").Value = -3254.6365
$ws.Range("N40").Value = -9271.909

$ws.Range("H63").Value = 28000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 28000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 28000
$ws.Range("N63").Value = -29498

$ws.Range("H66").Value = 28000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 28000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 84000
$ws.Range("N66").Value = -91488

$ws.Range("H93").Value = 1110.9
$ws.Range("I93").Value = 1012.7778
$ws.Range("J93").Value = 1994
$ws.Range("K93").Value = 1012.7778
$ws.Range("L93").Value = 1994
$ws.Range("M93").Value = 235.2222
$ws.Range("N93").Value = -4490

$ws.Range("H136").Value = 3201.639
$ws.Range("I136").Value = 2276.2964
$ws.Range("J136").Value = 5977.6665
$ws.Range("K136").Value = 6828.889200000001
$ws.Range("L136").Value = 17932.9995
$ws.Range("M136").Value = -4278.889200000001
$ws.Range("N136").Value = -23032.9995

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H81").Value = 219.8
$ws.Range("I81").Value = 199.75
$ws.Range("J81").Value = 300
$ws.Range("K81").Value = 399.5
$ws.Range("L81").Value = 600
$ws.Range("M81").Value = 661.5

$ws.Range("H82").Value = 50301
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 50301
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 50301
$ws.Range("N82").Value = -51067

$ws.Range("H84").Value = 219.8
$ws.Range("I84").Value = 199.75
$ws.Range("J84").Value = 300
$ws.Range("K84").Value = 1997.5
$ws.Range("L84").Value = 3000
$ws.Range("M84").Value = 3306.5

$ws.Range("H85").Value = 50301
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 50301
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 50301
$ws.Range("N85").Value = -52953

$ws.Range("H123").Value = 47533
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 47533
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 47533
$ws.Range("N123").Value = -57333

$ws.Range("H136").Value = 2918.0571
$ws.Range("I136").Value = 2500.8696
$ws.Range("J136").Value = 3717.6667
$ws.Range("K136").Value = 7502.6088
$ws.Range("L136").Value = 11153.0001
$ws.Range("M136").Value = -4952.6088
$ws.Range("N136").Value = -16253.0001
